# Daily attendance processing - 2025-10-21 07:42:01
# Normalises the "Recorded By" (column G) author ordering so "System" is
# always listed first, refreshes the outstanding/pending session counters,
# widens the Students column, and flags the three not-yet-run sessions
# (B2D / B2E / B2F, session 18) as "Not Recorded" with the pink status style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Widen column I (Students, 9th column) from 10 to 14 ---
# NOTE: the stored OOXML <col width> runs ~0.8333 wider than the COM
# ColumnWidth (character-count) value for this workbook's default font, so
# back that offset out here to land on an exact stored width of 14.
$ws.Columns.Item(9).ColumnWidth = 13.166666666666666

# --- 2. Re-order "Recorded By" text so "System" leads the list ---
$gUpdates = @{
    2 = "System, backup@backdoor.com, system"
    3 = "System, dnasr281@gmail.com"
    4 = "System, backup@backdoor.com"
    5 = "System, backup@backdoor.com"
    6 = "System, dnasr281@gmail.com"
    8 = "System, backup@backdoor.com"
    10 = "System, dnasr281@gmail.com"
    12 = "System, dnasr281@gmail.com"
    13 = "System, dnasr281@gmail.com"
    14 = "System, dnasr281@gmail.com"
    15 = "System, dnasr281@gmail.com"
    18 = "System, dnasr281@gmail.com"
    19 = "System, dnasr281@gmail.com"
    20 = "System, dnasr281@gmail.com"
    21 = "System, dnasr281@gmail.com"
    22 = "System, dnasr281@gmail.com"
    29 = "System, backup@backdoor.com, system"
    30 = "System, dnasr281@gmail.com"
    31 = "System, backup@backdoor.com"
    32 = "System, backup@backdoor.com"
    33 = "System, dnasr281@gmail.com"
    35 = "System, backup@backdoor.com"
    37 = "System, dnasr281@gmail.com"
    39 = "System, dnasr281@gmail.com"
    40 = "System, dnasr281@gmail.com"
    41 = "System, dnasr281@gmail.com"
    42 = "System, dnasr281@gmail.com"
    45 = "System, dnasr281@gmail.com"
    46 = "System, dnasr281@gmail.com"
    47 = "System, dnasr281@gmail.com"
    48 = "System, dnasr281@gmail.com"
    49 = "System, dnasr281@gmail.com"
    56 = "System, backup@backdoor.com, system"
    57 = "System, dnasr281@gmail.com"
    58 = "System, backup@backdoor.com"
    59 = "System, backup@backdoor.com"
    60 = "System, dnasr281@gmail.com"
    62 = "System, backup@backdoor.com"
    64 = "System, dnasr281@gmail.com"
    66 = "System, dnasr281@gmail.com"
    67 = "System, dnasr281@gmail.com"
    68 = "System, dnasr281@gmail.com"
    69 = "System, dnasr281@gmail.com"
    72 = "System, dnasr281@gmail.com"
    73 = "System, dnasr281@gmail.com"
    74 = "System, dnasr281@gmail.com"
    75 = "System, dnasr281@gmail.com"
    76 = "System, dnasr281@gmail.com"
    83 = "System, backup@backdoor.com"
    84 = "System, backup@backdoor.com"
    85 = "System, backup@backdoor.com"
    86 = "System, dnasr281@gmail.com"
    87 = "System, dnasr281@gmail.com"
    88 = "System, dnasr281@gmail.com"
    89 = "System, dnasr281@gmail.com"
    93 = "System, dnasr281@gmail.com"
    95 = "System, dnasr281@gmail.com"
    99 = "System, dnasr281@gmail.com"
    109 = "System, backup@backdoor.com"
    110 = "System, backup@backdoor.com"
    111 = "System, backup@backdoor.com"
    112 = "System, dnasr281@gmail.com"
    113 = "System, dnasr281@gmail.com"
    114 = "System, dnasr281@gmail.com"
    115 = "System, dnasr281@gmail.com"
    119 = "System, dnasr281@gmail.com"
    121 = "System, dnasr281@gmail.com"
    125 = "System, dnasr281@gmail.com"
    135 = "System, backup@backdoor.com"
    136 = "System, backup@backdoor.com"
    137 = "System, backup@backdoor.com"
    138 = "System, dnasr281@gmail.com"
    139 = "System, dnasr281@gmail.com"
    140 = "System, dnasr281@gmail.com"
    141 = "System, dnasr281@gmail.com"
    145 = "System, dnasr281@gmail.com"
    147 = "System, dnasr281@gmail.com"
    151 = "System, dnasr281@gmail.com"
}

foreach ($row in $gUpdates.Keys) {
    $ws.Range("G$row").Value = $gUpdates[$row]
}

# --- 3. Refresh the "Missing Sessions" / "Pending Sessions" counters ---
$ws.Range("L7").Value = 3
$ws.Range("L8").Value = 42

# --- 4. Refresh per-group "Pending"/"Not Recorded" breakdown (rows 18-20) ---
$ws.Range("P18").Value = 1
$ws.Range("Q18").Value = 8
$ws.Range("P19").Value = 1
$ws.Range("Q19").Value = 8
$ws.Range("P20").Value = 1
$ws.Range("Q20").Value = 8

# --- 5. Mark the 3 sessions that have not run yet as "Not Recorded" and
#        restyle their rows (A:H) with the pink "not-recorded" status style ---
$notRecordedRows = @(100, 126, 152)
foreach ($row in $notRecordedRows) {
    $rowRange = $ws.Range("A" + $row + ":H" + $row)
    $rowRange.Interior.Color = 12695295
    $rowRange.Font.Color = 0
    $ws.Range("I$row").Value = "Not Recorded"
    $ws.Range("I$row").Interior.Color = 12695295
    $ws.Range("I$row").Font.Color = 0
}
